$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8246
$ws1.Range("F5").Value = 6017
$ws1.Range("F6").Value = 515
$ws1.Range("F11").Value = 895

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8246
$ws4.Range("F5").Value = 6017
$ws4.Range("F6").Value = 515
$ws4.Range("F15").Value = 895
